$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C previously held the shared string "Home" for every data row;
# the reading-file fix relabels that column "Island" for each row (1-9).
$ws.Range("C1:C9").Value = "Island"

# Reflect the cursor position left behind after the edit.
$ws.Range("C10").Select()
